# Plantilla Lista de Tareas de la 7ma Iteración - update task status for
# "CU Consultar grupos y rentas" (row 7) to "En proceso", register 1 hour
# consumed on day 1 (K7), and leave the selection on K9 as the author did
# while reviewing the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 6 ("CU Generar recibo de pago.") status cell also moved from
# "Por iniciar" to the newly introduced "En proceso" status value.
$ws.Range("F6").Value = "En proceso"

# Row 7 ("CU Consultar grupos y rentas.") status + hours consumed on day 1.
$ws.Range("F7").Value = "En proceso"
$ws.Range("K7").Value = 1

# Re-create the header merges so the merge list is rebuilt with the same
# ranges the author ended up touching listed first (matches the saved
# workbook's mergeCells ordering).
$mergeOrder = @(
    "AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4",
    "AL4:AM4", "H4:I4", "K4:L4", "N4:O4", "Q4:R4", "T4:U4",
    "W4:X4", "Z4:AA4", "AC4:AD4", "AF4:AG4", "AI4:AJ4"
)

foreach ($addr in $mergeOrder) {
    $ws.Range($addr).UnMerge()
}
foreach ($addr in $mergeOrder) {
    $ws.Range($addr).Merge()
}

# Leave the cursor where the author last clicked.
$ws.Range("K9").Select()
